# Convert Short Partner Name to Upper case for proper corelation with BPA details table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPADetails")

# Column A (Partner) holds the short partner name used to correlate with the
# BPA details table. Normalise every value (rows 2-17) to upper case.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 17 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val.ToString().ToUpper()
    }
}

# --- Leftover UI/selection state from the editing session ---
# BPADetails ends up with F6 selected.
$ws.Activate()
$ws.Range("F6").Select()

# CmQuote had its columns selected (e.g. to inspect/clear formatting).
$wsCmQuote = $wb.Worksheets.Item("CmQuote")
$wsCmQuote.Activate()
$wsCmQuote.Range("A1:XFD1048576").Select()

# Shub ends up as the active/selected sheet with all columns selected.
$wsShub = $wb.Worksheets.Item("Shub")
$wsShub.Activate()
$wsShub.Range("A1:XFD1048576").Select()
